$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: replace formula-driven date in C22 with a static value, and add a new static date in D22
$ws.Range("C22").Value = 44014
$ws.Range("D22").Value = 44015
$ws.Range("D22").NumberFormat = $ws.Range("C22").NumberFormat

# Row 23: new label (" ") in C23, and a day-count formula in D23
$ws.Range("C23").Value = " "
$ws.Range("D23").Formula = "=D22-D12"
$ws.Range("D23").Style = "Normal"

# Row 24: update existing formulas to the new offsets
$ws.Range("C24").Formula = "=C13+98"
$ws.Range("D24").Formula = "=D13+100"

# Update the active selection to match the post-edit cursor position
$ws.Range("C25").Select()
